$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Theme Party" activity to "Kolokium Zon Selatan"
$ws.Range("C40").Value = "Kolokium Zon Selatan"

# Update figures
$ws.Range("E30").Value = 3500
$ws.Range("D40").Value = 2000

# Merge B15:C15 (new merged label cell for "Penandaan Fail" row)
$ws.Range("B15:C15").Merge()

# Match C15's formatting to the other plain bordered / non-wrapping cells
# in the table (e.g. C21, which already carries that border style) by
# copying the format only, after the merge so the merge state is preserved.
$ws.Range("C21").Copy()
$ws.Range("C15").PasteSpecial(-4122)
